# AlojamientoSeleccionado.xlsx - update listing name, host name, and price;
# widen column A to fit the new (longer) listing name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "I Loft You - Hermoso apartamento en Medellín"
$ws.Range("B2").Value = "Carolina"
$ws.Range("C2").Value = "$1,221,457.68 COP"

# Resize column A (bestFit) to accommodate the new, longer listing name.
$ws.Columns.Item(1).ColumnWidth = 41.8
